$wb = $excel.ActiveWorkbook
$wsMeas = $wb.Worksheets.Item("Measurements")
$wsStats = $wb.Worksheets.Item("Statistics")

# --- Update raw measurement values (Measurements sheet, rows 2-3) ---
$wsMeas.Range("A2").Value = 4.295692920684814
$wsMeas.Range("AA2").Value = 4.969364166
$wsMeas.Range("AB2").Value = -50.12388611
$wsMeas.Range("AC2").Value = 0.3214728832244873
$wsMeas.Range("AD2").Value = 5.01970291138
$wsMeas.Range("AE2").Value = -56.5685920715
$wsMeas.Range("AF2").Value = -56.6907081604
$wsMeas.Range("AG2").Value = 0.1915898323059082
$wsMeas.Range("AH2").Value = 10.32633566856384
$wsMeas.Range("AJ2").Value = 1.089
$wsMeas.Range("AL2").Value = 4.983222961
$wsMeas.Range("AM2").Value = -50.0286026
$wsMeas.Range("AN2").Value = 0.3172593116760254
$wsMeas.Range("AO2").Value = 5.03259134293
$wsMeas.Range("AP2").Value = -60.1378731728
$wsMeas.Range("AQ2").Value = -60.3335328102
$wsMeas.Range("AR2").Value = 0.1960489749908447
$wsMeas.Range("AS2").Value = 20.2311418056488
$wsMeas.Range("AU2").Value = 1.252
$wsMeas.Range("AW2").Value = 36.85
$wsMeas.Range("B2").Value = 12.95009994506836
$wsMeas.Range("F2").Value = 0.393
$wsMeas.Range("G2").Value = -13.4091
$wsMeas.Range("H2").Value = 5.068114
$wsMeas.Range("I2").Value = 4.973526001
$wsMeas.Range("J2").Value = -47.02700806
$wsMeas.Range("K2").Value = 0.1871702671051025
$wsMeas.Range("L2").Value = 5.01795244217
$wsMeas.Range("M2").Value = -54.9861989021
$wsMeas.Range("N2").Value = -56.114086628
$wsMeas.Range("O2").Value = 0.2370998859405518
$wsMeas.Range("P2").Value = 4.902143478
$wsMeas.Range("Q2").Value = -50.4039917
$wsMeas.Range("R2").Value = 0.1929292678833008
$wsMeas.Range("S2").Value = 4.96074104309
$wsMeas.Range("T2").Value = -61.3053684235
$wsMeas.Range("U2").Value = -60.9255847931
$wsMeas.Range("V2").Value = 0.1852684020996094
$wsMeas.Range("W2").Value = 5.775356769561768
$wsMeas.Range("Y2").Value = 1.563
$wsMeas.Range("A3").Value = 4.295692920684814
$wsMeas.Range("AA3").Value = 4.929424286
$wsMeas.Range("AB3").Value = -49.52594757
$wsMeas.Range("AC3").Value = 0.3138315677642822
$wsMeas.Range("AD3").Value = 5.01384782791
$wsMeas.Range("AE3").Value = -55.5860295296
$wsMeas.Range("AF3").Value = -55.7514500618
$wsMeas.Range("AG3").Value = 0.1920692920684814
$wsMeas.Range("AH3").Value = 10.48429465293884
$wsMeas.Range("AJ3").Value = 1.214
$wsMeas.Range("AL3").Value = 4.940135956
$wsMeas.Range("AM3").Value = -50.07219696
$wsMeas.Range("AN3").Value = 0.3605415821075439
$wsMeas.Range("AO3").Value = 5.03083229065
$wsMeas.Range("AP3").Value = -60.4961071014
$wsMeas.Range("AQ3").Value = -61.4769763947
$wsMeas.Range("AR3").Value = 0.2140405178070068
$wsMeas.Range("AS3").Value = 19.07639932632446
$wsMeas.Range("AU3").Value = 1.128
$wsMeas.Range("AW3").Value = 33.095
$wsMeas.Range("B3").Value = 12.95009994506836
$wsMeas.Range("C3").Value = 3
$wsMeas.Range("F3").Value = 0.352
$wsMeas.Range("G3").Value = -12.80008
$wsMeas.Range("H3").Value = 5.024055
$wsMeas.Range("I3").Value = 4.927913666
$wsMeas.Range("J3").Value = -50.06105423
$wsMeas.Range("K3").Value = 0.1960487365722656
$wsMeas.Range("L3").Value = 5.00654697418
$wsMeas.Range("M3").Value = -60.3990488052
$wsMeas.Range("N3").Value = -61.2072763443
$wsMeas.Range("O3").Value = 0.2091443538665771
$wsMeas.Range("P3").Value = 4.925983429
$wsMeas.Range("Q3").Value = -50.04518509
$wsMeas.Range("R3").Value = 0.1975250244140625
$wsMeas.Range("S3").Value = 5.01486253738
$wsMeas.Range("T3").Value = -60.3891682625
$wsMeas.Range("U3").Value = -61.2270207405
$wsMeas.Range("V3").Value = 0.1903145313262939
$wsMeas.Range("W3").Value = 3.019399642944336
$wsMeas.Range("Y3").Value = 1.28

# --- Update derived Statistics sheet (Max/Min/Mean per column) ---
$wsStats.Range("B3").Value = 4.295692920684814
$wsStats.Range("B4").Value = 4.295692920684814
$wsStats.Range("B5").Value = 4.295692920684814
$wsStats.Range("B6").Value = 12.95009994506836
$wsStats.Range("B7").Value = 12.95009994506836
$wsStats.Range("B8").Value = 12.95009994506836
$wsStats.Range("B9").Value = 3
$wsStats.Range("B10").Value = 2
$wsStats.Range("B11").Value = 2.5
$wsStats.Range("B18").Value = 0.393
$wsStats.Range("B19").Value = 0.352
$wsStats.Range("B20").Value = 0.3725
$wsStats.Range("B21").Value = -12.80008
$wsStats.Range("B22").Value = -13.4091
$wsStats.Range("B23").Value = -13.10459
$wsStats.Range("B24").Value = 5.068114
$wsStats.Range("B25").Value = 5.024055
$wsStats.Range("B26").Value = 5.046084499999999
$wsStats.Range("B27").Value = 4.973526001
$wsStats.Range("B28").Value = 4.927913666
$wsStats.Range("B29").Value = 4.9507198335
$wsStats.Range("B30").Value = -47.02700806
$wsStats.Range("B31").Value = -50.06105423
$wsStats.Range("B32").Value = -48.54403114500001
$wsStats.Range("B33").Value = 0.1960487365722656
$wsStats.Range("B34").Value = 0.1871702671051025
$wsStats.Range("B35").Value = 0.1916095018386841
$wsStats.Range("B36").Value = 5.01795244217
$wsStats.Range("B37").Value = 5.00654697418
$wsStats.Range("B38").Value = 5.012249708175
$wsStats.Range("B39").Value = -54.9861989021
$wsStats.Range("B40").Value = -60.3990488052
$wsStats.Range("B41").Value = -57.69262385365001
$wsStats.Range("B42").Value = -56.114086628
$wsStats.Range("B43").Value = -61.2072763443
$wsStats.Range("B44").Value = -58.66068148615
$wsStats.Range("B45").Value = 0.2370998859405518
$wsStats.Range("B46").Value = 0.2091443538665771
$wsStats.Range("B47").Value = 0.2231221199035645
$wsStats.Range("B48").Value = 4.925983429
$wsStats.Range("B49").Value = 4.902143478
$wsStats.Range("B50").Value = 4.914063453500001
$wsStats.Range("B51").Value = -50.04518509
$wsStats.Range("B52").Value = -50.4039917
$wsStats.Range("B53").Value = -50.224588395
$wsStats.Range("B54").Value = 0.1975250244140625
$wsStats.Range("B55").Value = 0.1929292678833008
$wsStats.Range("B56").Value = 0.1952271461486816
$wsStats.Range("B57").Value = 5.01486253738
$wsStats.Range("B58").Value = 4.96074104309
$wsStats.Range("B59").Value = 4.987801790235
$wsStats.Range("B60").Value = -60.3891682625
$wsStats.Range("B61").Value = -61.3053684235
$wsStats.Range("B62").Value = -60.847268343
$wsStats.Range("B63").Value = -60.9255847931
$wsStats.Range("B64").Value = -61.2270207405
$wsStats.Range("B65").Value = -61.0763027668
$wsStats.Range("B66").Value = 0.1903145313262939
$wsStats.Range("B67").Value = 0.1852684020996094
$wsStats.Range("B68").Value = 0.1877914667129517
$wsStats.Range("B69").Value = 5.775356769561768
$wsStats.Range("B70").Value = 3.019399642944336
$wsStats.Range("B71").Value = 4.397378206253052
$wsStats.Range("B75").Value = 1.563
$wsStats.Range("B76").Value = 1.28
$wsStats.Range("B77").Value = 1.4215
$wsStats.Range("B81").Value = 4.969364166
$wsStats.Range("B82").Value = 4.929424286
$wsStats.Range("B83").Value = 4.949394226
$wsStats.Range("B84").Value = -49.52594757
$wsStats.Range("B85").Value = -50.12388611
$wsStats.Range("B86").Value = -49.82491684
$wsStats.Range("B87").Value = 0.3214728832244873
$wsStats.Range("B88").Value = 0.3138315677642822
$wsStats.Range("B89").Value = 0.3176522254943848
$wsStats.Range("B90").Value = 5.01970291138
$wsStats.Range("B91").Value = 5.01384782791
$wsStats.Range("B92").Value = 5.016775369645
$wsStats.Range("B93").Value = -55.5860295296
$wsStats.Range("B94").Value = -56.5685920715
$wsStats.Range("B95").Value = -56.07731080054999
$wsStats.Range("B96").Value = -55.7514500618
$wsStats.Range("B97").Value = -56.6907081604
$wsStats.Range("B98").Value = -56.2210791111
$wsStats.Range("B99").Value = 0.1920692920684814
$wsStats.Range("B100").Value = 0.1915898323059082
$wsStats.Range("B101").Value = 0.1918295621871948
$wsStats.Range("B102").Value = 10.48429465293884
$wsStats.Range("B103").Value = 10.32633566856384
$wsStats.Range("B104").Value = 10.40531516075134
$wsStats.Range("B108").Value = 1.214
$wsStats.Range("B109").Value = 1.089
$wsStats.Range("B110").Value = 1.1515
$wsStats.Range("B114").Value = 4.983222961
$wsStats.Range("B115").Value = 4.940135956
$wsStats.Range("B116").Value = 4.9616794585
$wsStats.Range("B117").Value = -50.0286026
$wsStats.Range("B118").Value = -50.07219696
$wsStats.Range("B119").Value = -50.05039978
$wsStats.Range("B120").Value = 0.3605415821075439
$wsStats.Range("B121").Value = 0.3172593116760254
$wsStats.Range("B122").Value = 0.3389004468917847
$wsStats.Range("B123").Value = 5.03259134293
$wsStats.Range("B124").Value = 5.03083229065
$wsStats.Range("B125").Value = 5.031711816790001
$wsStats.Range("B126").Value = -60.1378731728
$wsStats.Range("B127").Value = -60.4961071014
$wsStats.Range("B128").Value = -60.3169901371
$wsStats.Range("B129").Value = -60.3335328102
$wsStats.Range("B130").Value = -61.4769763947
$wsStats.Range("B131").Value = -60.90525460245
$wsStats.Range("B132").Value = 0.2140405178070068
$wsStats.Range("B133").Value = 0.1960489749908447
$wsStats.Range("B134").Value = 0.2050447463989258
$wsStats.Range("B135").Value = 20.2311418056488
$wsStats.Range("B136").Value = 19.07639932632446
$wsStats.Range("B137").Value = 19.65377056598663
$wsStats.Range("B141").Value = 1.252
$wsStats.Range("B142").Value = 1.128
$wsStats.Range("B143").Value = 1.19
$wsStats.Range("B147").Value = 36.85
$wsStats.Range("B148").Value = 33.095
$wsStats.Range("B149").Value = 34.9725
